$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 rewrite: the "30 loka" diary entry is split into more granular
# columns — content note, quality note, code note and a new META remark —
# and the logged hours/time range are updated.

# C19: Oppimisen sisältö (learning content) — trimmed, page reference moved out
$ws.Range("C19").Value = "Initial draw ongelman selvittelyä fireworkscenestä"

# B19: Kello (time) — now covers two sessions, needs wrap text (style goes 5 -> 3)
$ws.Range("B19").Value = "9.15-11.15, 20.15-22.15"
$ws.Range("B19").WrapText = $true

# E19: Huomiot koodista (code notes) — new note
$ws.Range("E19").Value = "Vedetty muuttuvaa tilaa yhteen paikkaan ja oiottu skeneä hieman."
$ws.Range("E19").WrapText = $true

# D19: Oppimisen laatu (learning quality) — new note
$ws.Range("D19").Value = "noh tulipa sentään korjattua initial draw ongelmaa aika järeästi"
$ws.Range("D19").WrapText = $true

# F19: META — new remark (leading space preserved from source diff)
$ws.Range("F19").Value = " Kyllä tähän voisi melkein tottua että saa asioita aikaan"

# G19: Tunnit (hours) — 2 -> 4
$ws.Range("G19").Value = 4

# Row got shorter after trimming C19's text out to its own cell, now 2 wrapped lines
$ws.Rows.Item(19).RowHeight = 29

# Selection moves to F20, and the view no longer pins a frozen top-left cell
$ws.Range("F20").Select() | Out-Null
